$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Week of 12/6/2022 (column P) day-after inputs
$ws.Range("P3").Value = "DNP"
$ws.Range("P4").Value = "W"
$ws.Range("P5").Value = "W"
$ws.Range("P6").Value = "NA"
$ws.Range("P7").Value = "W"
$ws.Range("P8").Value = "W"
$ws.Range("P9").Value = "W"
$ws.Range("P10").Value = "DNP"

$ws.Range("P15").Value = "W"
$ws.Range("P16").Value = "DNP"
$ws.Range("P17").Value = "L"
$ws.Range("P18").Value = "DNP"
$ws.Range("P19").Value = "W"
$ws.Range("P20").Value = "NA"
$ws.Range("P21").Value = "DNP"
$ws.Range("P22").Value = "L"

# Update view state
$ws.Range("Q13").Select()
